$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "602.84") but must
# stay plain text, exactly like the original inline-string cells (this also
# preserves thousand-separator values such as "68.499.09" verbatim). Briefly
# force text format on each D cell right before writing it so Excel does not
# silently coerce the assigned string into a number, then restore the default
# "Normal" style so the cell formatting matches the untouched original.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.499.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.902.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.901.09"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  +0.56%  "
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.557.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.906.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.532.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.01%  "
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("E21").Value = "  -1.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "472.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.743"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000166"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.47%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("E30").Value = "  +1.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.052.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("E32").Value = "  +2.26%  "
$ws.Range("E33").Value = "  -2.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.49"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.874.29"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.41%  "
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.73"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +16.85%  "
$ws.Range("B39").Value = "Mantle"
$ws.Range("C39").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.03"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.142"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.314"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "428.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("B45").Value = "FLOKI"
$ws.Range("C45").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000302"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +12.77%  "
$ws.Range("E46").Value = "  +0.39%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.37%  "
$ws.Range("E49").Value = "  -2.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "27.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "143.85"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.73%  "
